$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New git "basic command" block (row 3) ---
$basicCommandBody = @'
$git rm {file}
$git status (-s)
$git diff
$git diff --cached
'@

# --- New "clone to local" block (row 4) ---
$cloneLocalBody = '$git clone https://github.com/mylesieong/my_maven_projects.git'

# --- New "clone to host" block (row 5) ---
$cloneHostBody = @'
**從無到有:
$git init
$git add {file}
$git config -global user.email "myles.ieong@gmail.com"
$git config -global user.name "myles"
$git commit -m "a project name"
$git remote add origin https://github.com/mylesieong/my_maven_projects.git
$git push -u origin master
(mylesieong:sewshort)
**從有到有:
$git push origin master
'@

# Write the new rows in the same order the strings were first introduced
# (bottom row up, column C before B before A) so the shared-string table
# grows in the same sequence as the authored workbook.

# Row 5
$ws.Range("C5").Value = $cloneHostBody
$ws.Range("B5").Value = "clone to host"
$ws.Range("A5").Value = "Git"

# Row 4
$ws.Range("C4").Value = $cloneLocalBody
$ws.Range("B4").Value = "clone to local"
$ws.Range("A4").Value = "Git"

# Row 3
$ws.Range("C3").Value = $basicCommandBody
$ws.Range("B3").Value = "basic command"
$ws.Range("A3").Value = "Git"

# Formatting for the new rows - same font family (Arial 10) used across A3:C5.
# Build the font once on a scratch cell and paste the format across the new
# range so a single new style entry is created instead of one per property.
$scratch = $ws.Range("F1")
$scratch.Font.Name = "Arial"
$scratch.Font.Size = 10
[void]$scratch.Copy()
$newRange = $ws.Range("A3:C5")
[void]$newRange.PasteSpecial(-4122)
[void]$scratch.Clear()

# Reflect the new used range / selection like the authored workbook
[void]$ws.Range("C7").Select()
